$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (originally sitting right after
#    "CAM: Camera").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. The paragraph "$TGT<SubPacketType>,L<Message>" is stored as two
#    adjacent runs with identical formatting. Re-run Find & Replace over it
#    (replacing the text with itself) so the engine regenerates it as a
#    single run while keeping the existing character formatting (color).
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$rng.Find.Execute("`$TGT<SubPacketType>,L<Message>", $true, $false, $false, $false, $false, $true, 1, $false, `
                   "`$TGT<SubPacketType>,L<Message>", 2) | Out-Null

# 3. Re-locate that text and drop a fresh "_GoBack" bookmark in the middle
#    of it, right after "$TGT<SubP" (splitting the run there), matching
#    "$TGT<SubP" + bookmark + "acketType>,L<Message>".
$rng2 = $d.Content
$rng2.Find.Execute("`$TGT<SubPacketType>,L<Message>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $rng2.Start + 9  # length of "$TGT<SubP"
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
